# Branched from CRAN release
# Adds a new "2.5.0" row (row 27) of benchmark results to the httk
# benchmarks table on Sheet1, grows Table1 to include it, and leaves the
# active selection on the new Notes cell - mirroring the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New data row (row 27) -------------------------------------------------
$ws.Range("A27").Value = "2.5.0"
$ws.Range("B27").Value = 1021
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0.9999
$ws.Range("F27").Value = 0.9477
$ws.Range("G27").Value = 353
$ws.Range("H27").Value = 0.2716
$ws.Range("I27").Value = 353
$ws.Range("J27").Value = 1.508
$ws.Range("K27").Value = 36
$ws.Range("L27").Value = 0.9698
$ws.Range("M27").Value = 80
$ws.Range("N27").Value = 1.132
$ws.Range("O27").Value = 80
$ws.Range("P27").Value = 0.6466
$ws.Range("Q27").Value = 863
$ws.Range("R27").Value = "Added models 3comp2 and sumclearances"

# --- Grow the Excel Table (ListObject) to include the new row -------------
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:R27"))

# --- Match the author's final selection on save ----------------------------
$ws.Range("R27").Select()
